# Test Suite Web - 14/05/2022
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Update the registered test-account text:
#    "halosalsa"            -> "halosalsa1"             (cell B7, the display text of a mailto hyperlink)
#    "halosalsa@gmail.com"  -> "halosalsa1@gmail.com"    (cell B8, the display text of a mailto hyperlink)
$ws.Range("B7").Value = "halosalsa1"
$ws.Range("B8").Value = "halosalsa1@gmail.com"

# 2. Widen column B (email column) so the longer address still fits comfortably.
$ws.Columns.Item(2).ColumnWidth = 29.29

# 3. Give the whole data range a neat grid: thin border around every cell plus
#    centered (horizontal + vertical) text alignment.
$used = $ws.Range("A1:D8")
$used.HorizontalAlignment = -4108   # xlCenter
$used.VerticalAlignment = -4108     # xlCenter
$used.Borders.LineStyle = 1         # xlContinuous
$used.Borders.Weight = 2            # xlThin

# 4. Highlight the header row with an orange fill.
$ws.Range("A1:D1").Interior.Color = 49407   # RGB(255,192,0) -> BGR long value
